$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report-date references: "January 2017" -> "February 2017" ---
# Main title (A2) and the chart's linked "Source:" caption cell (A148)
# both hold shared strings referencing the report month; update via the
# cells so the shared-string table (and the chart textbox field that
# links to A148) follow along.
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("A148").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Refreshed historical data for the most recent ~2 years (STEO Feb 2017 vintage) ---
$ws.Range("B122").Value = 172.13907900000001
$ws.Range("B123").Value = 170.00460000000001
$ws.Range("B124").Value = 168.32679999999999
$ws.Range("B125").Value = 168.2518
$ws.Range("B126").Value = 175.32650000000001
$ws.Range("B127").Value = 175.64930000000001
$ws.Range("B128").Value = 176.6497
$ws.Range("B129").Value = 170.07380000000001
$ws.Range("B130").Value = 160.655
$ws.Range("B131").Value = 154.85669999999999
$ws.Range("B132").Value = 152.56030000000001
$ws.Range("B133").Value = 156.22819999999999
$ws.Range("B134").Value = 160.25399999999999
$ws.Range("B135").Value = 156.31299999999999
$ws.Range("B136").Value = 150.53620000000001
$ws.Range("B137").Value = 149.57079999999999
$ws.Range("B138").Value = 156.16820000000001
$ws.Range("B139").Value = 156.9675
$ws.Range("B140").Value = 158.2662
$ws.Range("B141").Value = 152.9802
$ws.Range("B142").Value = 144.84289999999999
$ws.Range("B143").Value = 139.81909999999999
$ws.Range("B144").Value = 137.38900000000001
$ws.Range("B145").Value = 141.71530000000001
$ws.Range("B146").Value = 146.49299999999999
$ws.Range("B147").Value = 148.7978

# --- Forecast marker (scatter series x-values): shift from month-index 96 to 97 ---
$ws.Range("A152").Value = 97
$ws.Range("A153").Value = 97
